$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data occupies rows 2-28 in columns A and B.
# Fix the iteration count: shift every row's values up by one row
# (each row takes on the values that used to belong to the next row),
# leaving the final row (28) unchanged.
$firstRow = 2
$lastRow = 28

for ($r = $firstRow; $r -lt $lastRow; $r++) {
    $nextA = $ws.Cells.Item($r + 1, 1).Value2
    $nextB = $ws.Cells.Item($r + 1, 2).Value2
    $ws.Cells.Item($r, 1).Value = $nextA
    $ws.Cells.Item($r, 2).Value = $nextB
}
